# Apply updated TPM-derived values to C3-Lrp1 LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"3.820425"
$ws.Range("H2").Value = [double]"11.461275"
$ws.Range("I2").Value = [double]"0.02049663039797357"
$ws.Range("J2").Value = [double]"0.02049663039797357"
$ws.Range("M2").Value = [double]"6.305846"
$ws.Range("N2").Value = [double]"18.917538"
$ws.Range("O2").Value = [double]"0.01356150511917599"
$ws.Range("P2").Value = [double]"0.01356150511917599"
$ws.Range("Q2").Value = [double]"24.09101170455"
$ws.Range("R2").Value = [double]"216.81910534095"
$ws.Range("S2").Value = [double]"0.0002779651580679767"
$ws.Range("T2").Value = [double]"0.0002779651580679766"
$ws.Range("G3").Value = [double]"3.820425"
$ws.Range("H3").Value = [double]"11.461275"
$ws.Range("I3").Value = [double]"0.02049663039797357"
$ws.Range("J3").Value = [double]"0.02049663039797357"
$ws.Range("O3").Value = [double]"0.392557056479861"
$ws.Range("P3").Value = [double]"0.3925570564798609"
$ws.Range("Q3").Value = [double]"697.3486024783251"
$ws.Range("R3").Value = [double]"6276.137422304925"
$ws.Range("S3").Value = [double]"0.008046096896784144"
$ws.Range("T3").Value = [double]"0.008046096896784144"
$ws.Range("G4").Value = [double]"3.820425"
$ws.Range("H4").Value = [double]"11.461275"
$ws.Range("I4").Value = [double]"0.02049663039797357"
$ws.Range("J4").Value = [double]"0.02049663039797357"
$ws.Range("M4").Value = [double]"127.396393"
$ws.Range("N4").Value = [double]"382.189179"
$ws.Range("O4").Value = [double]"0.2739817680029065"
$ws.Range("P4").Value = [double]"0.2739817680029065"
$ws.Range("Q4").Value = [double]"486.708364727025"
$ws.Range("R4").Value = [double]"4380.375282543226"
$ws.Range("S4").Value = [double]"0.005615703034538915"
$ws.Range("T4").Value = [double]"0.005615703034538915"
$ws.Range("G5").Value = [double]"3.820425"
$ws.Range("H5").Value = [double]"11.461275"
$ws.Range("I5").Value = [double]"0.02049663039797357"
$ws.Range("J5").Value = [double]"0.02049663039797357"
$ws.Range("M5").Value = [double]"19.42400133333333"
$ws.Range("N5").Value = [double]"58.272004"
$ws.Range("O5").Value = [double]"0.04177372766745037"
$ws.Range("P5").Value = [double]"0.04177372766745036"
$ws.Range("Q5").Value = [double]"74.20794029390001"
$ws.Range("R5").Value = [double]"667.8714626451001"
$ws.Range("S5").Value = [double]"0.0008562206563453326"
$ws.Range("T5").Value = [double]"0.0008562206563453325"
$ws.Range("G6").Value = [double]"3.820425"
$ws.Range("H6").Value = [double]"11.461275"
$ws.Range("I6").Value = [double]"0.02049663039797357"
$ws.Range("J6").Value = [double]"0.02049663039797357"
$ws.Range("M6").Value = [double]"129.3233566666667"
$ws.Range("N6").Value = [double]"387.97007"
$ws.Range("O6").Value = [double]"0.2781259427306063"
$ws.Range("P6").Value = [double]"0.2781259427306062"
$ws.Range("Q6").Value = [double]"494.0701848932501"
$ws.Range("R6").Value = [double]"4446.631664039251"
$ws.Range("S6").Value = [double]"0.0057006446522372"
$ws.Range("T6").Value = [double]"0.005700644652237198"
$ws.Range("I7").Value = [double]"0.7542622677884155"
$ws.Range("J7").Value = [double]"0.7542622677884157"
$ws.Range("M7").Value = [double]"6.305846"
$ws.Range("N7").Value = [double]"18.917538"
$ws.Range("O7").Value = [double]"0.01356150511917599"
$ws.Range("P7").Value = [double]"0.01356150511917599"
$ws.Range("Q7").Value = [double]"886.5330919655772"
$ws.Range("R7").Value = [double]"7978.797827690197"
$ws.Range("S7").Value = [double]"0.01022893160581389"
$ws.Range("T7").Value = [double]"0.01022893160581389"
$ws.Range("I8").Value = [double]"0.7542622677884155"
$ws.Range("J8").Value = [double]"0.7542622677884157"
$ws.Range("O8").Value = [double]"0.392557056479861"
$ws.Range("P8").Value = [double]"0.3925570564798609"
$ws.Range("S8").Value = [double]"0.2960909756568451"
$ws.Range("T8").Value = [double]"0.2960909756568451"
$ws.Range("I9").Value = [double]"0.7542622677884155"
$ws.Range("J9").Value = [double]"0.7542622677884157"
$ws.Range("M9").Value = [double]"127.396393"
$ws.Range("N9").Value = [double]"382.189179"
$ws.Range("O9").Value = [double]"0.2739817680029065"
$ws.Range("P9").Value = [double]"0.2739817680029065"
$ws.Range("Q9").Value = [double]"17910.54177211937"
$ws.Range("R9").Value = [double]"161194.8759490743"
$ws.Range("S9").Value = [double]"0.2066541096665518"
$ws.Range("T9").Value = [double]"0.2066541096665518"
$ws.Range("I10").Value = [double]"0.7542622677884155"
$ws.Range("J10").Value = [double]"0.7542622677884157"
$ws.Range("M10").Value = [double]"19.42400133333333"
$ws.Range("N10").Value = [double]"58.272004"
$ws.Range("O10").Value = [double]"0.04177372766745037"
$ws.Range("P10").Value = [double]"0.04177372766745036"
$ws.Range("Q10").Value = [double]"2730.802490321441"
$ws.Range("R10").Value = [double]"24577.22241289297"
$ws.Range("S10").Value = [double]"0.0315083465644268"
$ws.Range("T10").Value = [double]"0.0315083465644268"
$ws.Range("I11").Value = [double]"0.7542622677884155"
$ws.Range("J11").Value = [double]"0.7542622677884157"
$ws.Range("M11").Value = [double]"129.3233566666667"
$ws.Range("N11").Value = [double]"387.97007"
$ws.Range("O11").Value = [double]"0.2781259427306063"
$ws.Range("P11").Value = [double]"0.2781259427306062"
$ws.Range("Q11").Value = [double]"18181.45182249411"
$ws.Range("R11").Value = [double]"163633.0664024469"
$ws.Range("S11").Value = [double]"0.2097799042947781"
$ws.Range("T11").Value = [double]"0.2097799042947781"
$ws.Range("G12").Value = [double]"30.51067"
$ws.Range("H12").Value = [double]"91.53201"
$ws.Range("I12").Value = [double]"0.1636901460399144"
$ws.Range("J12").Value = [double]"0.1636901460399144"
$ws.Range("M12").Value = [double]"6.305846"
$ws.Range("N12").Value = [double]"18.917538"
$ws.Range("O12").Value = [double]"0.01356150511917599"
$ws.Range("P12").Value = [double]"0.01356150511917599"
$ws.Range("Q12").Value = [double]"192.39558637682"
$ws.Range("R12").Value = [double]"1731.56027739138"
$ws.Range("S12").Value = [double]"0.002219884753478965"
$ws.Range("T12").Value = [double]"0.002219884753478965"
$ws.Range("G13").Value = [double]"30.51067"
$ws.Range("H13").Value = [double]"91.53201"
$ws.Range("I13").Value = [double]"0.1636901460399144"
$ws.Range("J13").Value = [double]"0.1636901460399144"
$ws.Range("O13").Value = [double]"0.392557056479861"
$ws.Range("P13").Value = [double]"0.3925570564798609"
$ws.Range("Q13").Value = [double]"5569.163924217164"
$ws.Range("R13").Value = [double]"50122.47531795447"
$ws.Range("S13").Value = [double]"0.0642577219041874"
$ws.Range("T13").Value = [double]"0.06425772190418738"
$ws.Range("G14").Value = [double]"30.51067"
$ws.Range("H14").Value = [double]"91.53201"
$ws.Range("I14").Value = [double]"0.1636901460399144"
$ws.Range("J14").Value = [double]"0.1636901460399144"
$ws.Range("M14").Value = [double]"127.396393"
$ws.Range("N14").Value = [double]"382.189179"
$ws.Range("O14").Value = [double]"0.2739817680029065"
$ws.Range("P14").Value = [double]"0.2739817680029065"
$ws.Range("Q14").Value = [double]"3886.94930601331"
$ws.Range("R14").Value = [double]"34982.54375411979"
$ws.Range("S14").Value = [double]"0.04484811561666972"
$ws.Range("T14").Value = [double]"0.04484811561666972"
$ws.Range("G15").Value = [double]"30.51067"
$ws.Range("H15").Value = [double]"91.53201"
$ws.Range("I15").Value = [double]"0.1636901460399144"
$ws.Range("J15").Value = [double]"0.1636901460399144"
$ws.Range("M15").Value = [double]"19.42400133333333"
$ws.Range("N15").Value = [double]"58.272004"
$ws.Range("O15").Value = [double]"0.04177372766745037"
$ws.Range("P15").Value = [double]"0.04177372766745036"
$ws.Range("Q15").Value = [double]"592.6392947608933"
$ws.Range("R15").Value = [double]"5333.75365284804"
$ws.Range("S15").Value = [double]"0.006837947582516566"
$ws.Range("T15").Value = [double]"0.006837947582516564"
$ws.Range("G16").Value = [double]"30.51067"
$ws.Range("H16").Value = [double]"91.53201"
$ws.Range("I16").Value = [double]"0.1636901460399144"
$ws.Range("J16").Value = [double]"0.1636901460399144"
$ws.Range("M16").Value = [double]"129.3233566666667"
$ws.Range("N16").Value = [double]"387.97007"
$ws.Range("O16").Value = [double]"0.2781259427306063"
$ws.Range("P16").Value = [double]"0.2781259427306062"
$ws.Range("Q16").Value = [double]"3945.742258548967"
$ws.Range("R16").Value = [double]"35511.6803269407"
$ws.Range("S16").Value = [double]"0.04552647618306182"
$ws.Range("T16").Value = [double]"0.04552647618306181"
$ws.Range("G17").Value = [double]"0.258813"
$ws.Range("H17").Value = [double]"0.776439"
$ws.Range("I17").Value = [double]"0.001388535150720334"
$ws.Range("J17").Value = [double]"0.001388535150720334"
$ws.Range("M17").Value = [double]"6.305846"
$ws.Range("N17").Value = [double]"18.917538"
$ws.Range("O17").Value = [double]"0.01356150511917599"
$ws.Range("P17").Value = [double]"0.01356150511917599"
$ws.Range("Q17").Value = [double]"1.632034920798"
$ws.Range("R17").Value = [double]"14.688314287182"
$ws.Range("S17").Value = [double]"1.883062655464962E-05"
$ws.Range("T17").Value = [double]"1.883062655464961E-05"
$ws.Range("G18").Value = [double]"0.258813"
$ws.Range("H18").Value = [double]"0.776439"
$ws.Range("I18").Value = [double]"0.001388535150720334"
$ws.Range("J18").Value = [double]"0.001388535150720334"
$ws.Range("O18").Value = [double]"0.392557056479861"
$ws.Range("P18").Value = [double]"0.3925570564798609"
$ws.Range("Q18").Value = [double]"47.24157229973701"
$ws.Range("R18").Value = [double]"425.174150697633"
$ws.Range("S18").Value = [double]"0.0005450792715855946"
$ws.Range("T18").Value = [double]"0.0005450792715855944"
$ws.Range("G19").Value = [double]"0.258813"
$ws.Range("H19").Value = [double]"0.776439"
$ws.Range("I19").Value = [double]"0.001388535150720334"
$ws.Range("J19").Value = [double]"0.001388535150720334"
$ws.Range("M19").Value = [double]"127.396393"
$ws.Range("N19").Value = [double]"382.189179"
$ws.Range("O19").Value = [double]"0.2739817680029065"
$ws.Range("P19").Value = [double]"0.2739817680029065"
$ws.Range("Q19").Value = [double]"32.971842661509"
$ws.Range("R19").Value = [double]"296.746583953581"
$ws.Range("S19").Value = [double]"0.0003804333155285394"
$ws.Range("T19").Value = [double]"0.0003804333155285394"
$ws.Range("G20").Value = [double]"0.258813"
$ws.Range("H20").Value = [double]"0.776439"
$ws.Range("I20").Value = [double]"0.001388535150720334"
$ws.Range("J20").Value = [double]"0.001388535150720334"
$ws.Range("M20").Value = [double]"19.42400133333333"
$ws.Range("N20").Value = [double]"58.272004"
$ws.Range("O20").Value = [double]"0.04177372766745037"
$ws.Range("P20").Value = [double]"0.04177372766745036"
$ws.Range("Q20").Value = [double]"5.027184057084"
$ws.Range("R20").Value = [double]"45.244656513756"
$ws.Range("S20").Value = [double]"5.80042892428734E-05"
$ws.Range("T20").Value = [double]"5.800428924287338E-05"
$ws.Range("G21").Value = [double]"0.258813"
$ws.Range("H21").Value = [double]"0.776439"
$ws.Range("I21").Value = [double]"0.001388535150720334"
$ws.Range("J21").Value = [double]"0.001388535150720334"
$ws.Range("M21").Value = [double]"129.3233566666667"
$ws.Range("N21").Value = [double]"387.97007"
$ws.Range("O21").Value = [double]"0.2781259427306063"
$ws.Range("P21").Value = [double]"0.2781259427306062"
$ws.Range("Q21").Value = [double]"33.47056590897001"
$ws.Range("R21").Value = [double]"301.23509318073"
$ws.Range("S21").Value = [double]"0.0003861876478086774"
$ws.Range("T21").Value = [double]"0.0003861876478086773"
$ws.Range("G22").Value = [double]"11.213844"
$ws.Range("H22").Value = [double]"33.641532"
$ws.Range("I22").Value = [double]"0.0601624206229761"
$ws.Range("J22").Value = [double]"0.0601624206229761"
$ws.Range("M22").Value = [double]"6.305846"
$ws.Range("N22").Value = [double]"18.917538"
$ws.Range("O22").Value = [double]"0.01356150511917599"
$ws.Range("P22").Value = [double]"0.01356150511917599"
$ws.Range("Q22").Value = [double]"70.712773332024"
$ws.Range("R22").Value = [double]"636.414959988216"
$ws.Range("S22").Value = [double]"0.0008158929752605094"
$ws.Range("T22").Value = [double]"0.0008158929752605093"
$ws.Range("G23").Value = [double]"11.213844"
$ws.Range("H23").Value = [double]"33.641532"
$ws.Range("I23").Value = [double]"0.0601624206229761"
$ws.Range("J23").Value = [double]"0.0601624206229761"
$ws.Range("O23").Value = [double]"0.392557056479861"
$ws.Range("P23").Value = [double]"0.3925570564798609"
$ws.Range("Q23").Value = [double]"2046.881810743556"
$ws.Range("R23").Value = [double]"18421.93629669201"
$ws.Range("S23").Value = [double]"0.02361718275045878"
$ws.Range("T23").Value = [double]"0.02361718275045878"
$ws.Range("G24").Value = [double]"11.213844"
$ws.Range("H24").Value = [double]"33.641532"
$ws.Range("I24").Value = [double]"0.0601624206229761"
$ws.Range("J24").Value = [double]"0.0601624206229761"
$ws.Range("M24").Value = [double]"127.396393"
$ws.Range("N24").Value = [double]"382.189179"
$ws.Range("O24").Value = [double]"0.2739817680029065"
$ws.Range("P24").Value = [double]"0.2739817680029065"
$ws.Range("Q24").Value = [double]"1428.603277264692"
$ws.Range("R24").Value = [double]"12857.42949538223"
$ws.Range("S24").Value = [double]"0.01648340636961752"
$ws.Range("T24").Value = [double]"0.01648340636961752"
$ws.Range("G25").Value = [double]"11.213844"
$ws.Range("H25").Value = [double]"33.641532"
$ws.Range("I25").Value = [double]"0.0601624206229761"
$ws.Range("J25").Value = [double]"0.0601624206229761"
$ws.Range("M25").Value = [double]"19.42400133333333"
$ws.Range("N25").Value = [double]"58.272004"
$ws.Range("O25").Value = [double]"0.04177372766745037"
$ws.Range("P25").Value = [double]"0.04177372766745036"
$ws.Range("Q25").Value = [double]"217.817720807792"
$ws.Range("R25").Value = [double]"1960.359487270128"
$ws.Range("S25").Value = [double]"0.002513208574918803"
$ws.Range("T25").Value = [double]"0.002513208574918803"
$ws.Range("G26").Value = [double]"11.213844"
$ws.Range("H26").Value = [double]"33.641532"
$ws.Range("I26").Value = [double]"0.0601624206229761"
$ws.Range("J26").Value = [double]"0.0601624206229761"
$ws.Range("M26").Value = [double]"129.3233566666667"
$ws.Range("N26").Value = [double]"387.97007"
$ws.Range("O26").Value = [double]"0.2781259427306063"
$ws.Range("P26").Value = [double]"0.2781259427306062"
$ws.Range("Q26").Value = [double]"1450.21194721636"
$ws.Range("R26").Value = [double]"13051.90752494724"
$ws.Range("S26").Value = [double]"0.0167327299527205"
$ws.Range("T26").Value = [double]"0.01673272995272049"

Write-Host "Updated $(278) cells with new TPM values"
